{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: \"En tant qu'\u00e9tudiant, vous disposez d'un compte utilisateur sur LOSPA\u00c9\"\n//        -> \"En tant qu'\u00e9l\u00e8ve, vous disposez d'un compte utilisateur sur LOSPA\u00c9\"\n// Change 2: \"Bonne continuation en BTS SIO !\"\n//        -> \"Bonne continuation !\"\n\nconst body = context.document.body;\n\n// --- Change 1: \"\u00e9tudiant\" -> \"\u00e9l\u00e8ve\" ---\nconst studentSearch = body.search(\"qu\\u2019\u00e9tudiant\", { matchCase: true, matchWholeWord: false });\nstudentSearch.load(\"text\");\nawait context.sync();\n\nif (studentSearch.items.length > 0) {\n  studentSearch.items[0].insertText(\"qu\\u2019\u00e9l\u00e8ve\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: drop \" en BTS SIO\" from the closing line ---\nconst closingSearch = body.search(\"Bonne continuation en BTS SIO !\", { matchCase: true, matchWholeWord: false });\nclosingSearch.load(\"text\");\nawait context.sync();\n\nif (closingSearch.items.length > 0) {\n  closingSearch.items[0].insertText(\"Bonne continuation\\u00a0!\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: \"En tant qu'\u00e9tudiant, vous disposez d'un compte utilisateur sur LOSPA\u00c9\"\n#        -> \"En tant qu'\u00e9l\u00e8ve, vous disposez d'un compte utilisateur sur LOSPA\u00c9\"\n# Change 2: \"Bonne continuation en BTS SIO !\"\n#        -> \"Bonne continuation !\"\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"\u00e9tudiant\" -> \"\u00e9l\u00e8ve\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$result1 = $find1.Execute(\"qu\u2019\u00e9tudiant\", $false, $false, $false, $false, $false, $true, 1, $null, \"qu\u2019\u00e9l\u00e8ve\", 2)\n\n# --- Change 2: drop \" en BTS SIO\" from the closing line ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$result2 = $find2.Execute(\"Bonne continuation en BTS SIO !\", $false, $false, $false, $false, $false, $true, 1, $null, \"Bonne continuation\u00a0!\", 2)\n\nWrite-Output \"replace1=$result1 replace2=$result2\"\n"}
